$d = $word.ActiveDocument

$pairs = @(
    @("392÷3=", "296÷7="),
    @("692÷2=", "401÷8="),
    @("313÷3=", "502÷5="),
    @("606÷8=", "293÷7="),
    @("504÷8=", "882÷8="),
    @("201÷2=", "159÷2="),
    @("502÷7=", "510÷6="),
    @("815÷4=", "655÷8="),
    @("721÷3=", "334÷7="),
    @("946÷6=", "497÷2="),
    @("598÷4=", "875÷3="),
    @("908÷2=", "613÷8="),
    @("871÷6=", "640÷4="),
    @("674÷6=", "114÷2="),
    @("320÷9=", "296÷4="),
    @("575÷6=", "640÷9="),
    @("563÷3=", "750÷4="),
    @("644÷4=", "327÷5="),
    @("482÷8=", "821÷8="),
    @("262÷3=", "284÷8="),
    @("331÷6=", "272÷5="),
    @("739÷5=", "696÷5="),
    @("492÷2=", "656÷6="),
    @("316÷2=", "531÷5="),
    @("463÷8=", "125÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
